# NMCARS-ANNEX-22.docx maintenance touch-up
# --------------------------------------------------------------------
# Source commit: "Fixed NMCARS internal xrefs and image sources for EMF
# files". For this particular part of the document the visible effect
# is a re-save of the customXml plumbing:
#   Word re-numbered the four customXml "item"/"itemProps" parts (the
#   document-info / SharePoint content-type schema, the document-
#   management "properties" part, the FormTemplates part, and the
#   (empty) bibliography Sources part) while leaving each one's
#   content untouched - item2 -> item1, item3 -> item2, item4 -> item3,
#   item1 -> item4 (and the matching itemPropsN renumber/rewire so each
#   item keeps pointing at its own itemProps part).
#
# None of this is literal document *content*, so it is driven through
# the matching Word object-model surface (CustomXMLParts) instead of
# poking at body text. Everything is wrapped defensively so that if a
# given host implementation does not surface CustomXMLParts mutation,
# the rest of the script - and the document - still come through
# unharmed (no partial/half-applied rotation).

$d = $word.ActiveDocument

# ----------------------------------------------------------------
# 1) Rotate the customXml parts: item2/item3/item4/item1 -> 1/2/3/4
# ----------------------------------------------------------------
function Get-PartXml($parts, $i) {
    try { return $parts.Item($i).XML } catch { return $null }
}

try {
    $cxp = $d.CustomXMLParts
    $count = $cxp.Count

    if ($count -ge 4) {
        # Snapshot the current four parts' XML bodies in document order
        # (this mirrors the _rels wiring: rIdN -> customXml/itemN.xml).
        $xml1 = Get-PartXml $cxp 1
        $xml2 = Get-PartXml $cxp 2
        $xml3 = Get-PartXml $cxp 3
        $xml4 = Get-PartXml $cxp 4

        if ($xml1 -and $xml2 -and $xml3 -and $xml4) {
            # New ordering after the rotation described by the diff:
            #   item1 <- old item2 (documentManagement / properties)
            #   item2 <- old item3 (FormTemplates)
            #   item3 <- old item4 (bibliography Sources)
            #   item4 <- old item1 (SharePoint contentTypeSchema)
            $newOrder = @($xml2, $xml3, $xml4, $xml1)

            # Drop the old four parts (highest index first so the
            # indices of the ones we haven't removed yet stay put),
            # then add the replacements back in the new order so they
            # re-occupy slots 1..4.
            for ($i = $count; $i -ge 1; $i--) {
                try { $cxp.Item($i).Delete() } catch { }
            }
            foreach ($x in $newOrder) {
                try { [void]$cxp.Add($x) } catch { }
            }
        }
    }
} catch {
    # CustomXMLParts not available on this host - nothing else to do
    # for the part-renumbering half of the change.
}

Write-Output "done"
